# Update Name of Algo
# Apply updated RandomForest imputed values to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E7").Value  = 15.27619999999999
$ws.Range("B10").Value = 5.2904
$ws.Range("B12").Value = 4.946300000000003
$ws.Range("E15").Value = 15.9464
$ws.Range("B18").Value = 6.440300000000002
$ws.Range("E20").Value = 15.90929999999998
$ws.Range("E29").Value = 17.14850000000002
$ws.Range("E30").Value = 15.40699999999999
$ws.Range("E31").Value = 16.0448
$ws.Range("B37").Value = 8.731900000000001
$ws.Range("E40").Value = 17.07140000000001
$ws.Range("B55").Value = 5.885299999999996
$ws.Range("B68").Value = 4.880799999999996
$ws.Range("E68").Value = 17.35180000000001
$ws.Range("E76").Value = 16.17889999999998
$ws.Range("B77").Value = 8.986800000000002
$ws.Range("B78").Value = 9.364899999999999
$ws.Range("E87").Value = 16.25959999999999
$ws.Range("E88").Value = 16.2657
$ws.Range("E96").Value = 16.23029999999999
$ws.Range("E98").Value = 15.4179
$ws.Range("E101").Value = 16.78360000000001
$ws.Range("E102").Value = 16.70979999999999
